# Modified example grading sheet
#
# - Grades the "Qualite du code" criterion (row 6) for student #2
#   (column F) and student #5 (column I): 2 pts and 4 pts respectively.
#   All dependent totals/averages/notes (rows 4, 20, 21, 28) recalculate
#   automatically from these two inputs, since they are plain SUM /
#   AVERAGE formulas over the grade columns.
# - Leaves grading feedback as cell comments: a critical note on F6
#   ("Not very good work. (-2 pts)") and a praising note on H6
#   ("Great work !"), both signed "Ylli:" as the first line.
# - Leaves the active selection on F6, matching the reviewer's last
#   focused cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grade entry ---------------------------------------------------------
$ws.Range("F6").Value = 2
$ws.Range("I6").Value = 4

# --- Reviewer comments -----------------------------------------------------
$commentF6 = $ws.Range("F6").AddComment()
$commentF6.Text("Ylli:" + [char]10 + "Not very good work. (-2 pts)") | Out-Null

$commentH6 = $ws.Range("H6").AddComment()
$commentH6.Text("Ylli:" + [char]10 + "Great work !") | Out-Null

# --- Selection ---------------------------------------------------------
$ws.Range("F6").Select() | Out-Null
